$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.711.93"
$ws.Range("E2").Value = "  -0.89%  "
$ws.Range("D3").Value = "3.419.15"
$ws.Range("E3").Value = "  -2.43%  "
$ws.Range("E4").Value = "  +0.04%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "579.89"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -1.62%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "128.99"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -3.94%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.480"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -1.48%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "7.57"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +3.69%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.124"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -0.01%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.382"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -1.03%  "
$ws.Range("D12").Value = "4.002.04"
$ws.Range("E12").Value = "  -2.41%  "
$ws.Range("E13").Value = "  -0.42%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.0000176"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -2.39%  "
$ws.Range("D15").Value = "3.420.37"
$ws.Range("E15").Value = "  -2.42%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "63.707.25"
$ws.Range("E16").Value = "  -0.94%  "
$ws.Range("B17").Value = "Avalanche"
$ws.Range("C17").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "25.44"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -0.98%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "9.82"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -0.36%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "5.65"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -1.85%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "13.34"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -1.41%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "383.11"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -2.55%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.564"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -1.23%  "
$ws.Range("D23").Value = "3.556.71"
$ws.Range("E23").Value = "  -2.39%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "74.03"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -0.84%  "
$ws.Range("E25").Value = "  +0.16%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "0.0000109"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -4.86%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +0.23%  "
$ws.Range("E28").Value = "  -2.41%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "7.01"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -4.62%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "7.88"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -4.22%  "
$ws.Range("E31").Value = "  -0.47%  "
$ws.Range("E32").Value = "  -3.61%  "
$ws.Range("D33").Value = "3.448.89"
$ws.Range("E33").Value = "  -2.20%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "22.78"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -2.89%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "5.13"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -0.11%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "6.73"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -2.16%  "
$ws.Range("B38").Value = "Monero"
$ws.Range("C38").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "164.04"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -2.01%  "
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "1.51"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -2.32%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.0769"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -1.41%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.785"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -3.24%  "
$ws.Range("E42").Value = "  +0.04%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "41.41"
$c.Style = "Normal"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "4.31"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -1.93%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "1.60"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -3.52%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "23.31"
$c.Style = "Normal"
$ws.Range("E46").Value = "  -7.73%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "1.09"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -6.11%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "6.70"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -0.88%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.888"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -0.41%  "
$ws.Range("D50").Value = "2.272.05"
$ws.Range("E50").Value = "  -2.72%  "
$ws.Range("E51").Value = "  -2.36%  "
